$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column width tweaks ---
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(11).ColumnWidth = 10.3

# --- Hide rows 34-41 (week-old trades, now out of filter window) ---
for ($r = 34; $r -le 41; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# --- Hide rows 46-47 and 49-63 (same reason) ---
$ws.Rows.Item(46).Hidden = $true
$ws.Rows.Item(47).Hidden = $true
for ($r = 49; $r -le 63; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# --- New trade rows 64-67: ID 1452-70146 (30-Apr-2021 expiry SPX straddle) ---
$ws.Cells.Item(64,1).Value = 44309
$ws.Cells.Item(64,2).Value = "1452-70146"
$ws.Cells.Item(64,3).Value = "SELL"
$ws.Cells.Item(64,4).Value = 100
$ws.Cells.Item(64,5).Value = "SPX"
$ws.Cells.Item(64,6).Value = 44316
$ws.Cells.Item(64,7).Value = "CALL"
$ws.Cells.Item(64,8).Value = 4170
$ws.Cells.Item(64,9).Value = 22.83
$ws.Cells.Item(64,11).Formula = "=D64*I64"

$ws.Cells.Item(65,1).Value = 44309
$ws.Cells.Item(65,2).Value = "1452-70146"
$ws.Cells.Item(65,3).Value = "SELL"
$ws.Cells.Item(65,4).Value = 100
$ws.Cells.Item(65,5).Value = "SPX"
$ws.Cells.Item(65,6).Value = 44316
$ws.Cells.Item(65,7).Value = "PUT"
$ws.Cells.Item(65,8).Value = 4140
$ws.Cells.Item(65,9).Value = 24.99
$ws.Cells.Item(65,11).Formula = "=D65*I65"

$ws.Cells.Item(66,1).Value = 44309
$ws.Cells.Item(66,2).Value = "1452-70146"
$ws.Cells.Item(66,3).Value = "BUY"
$ws.Cells.Item(66,4).Value = -100
$ws.Cells.Item(66,5).Value = "SPX"
$ws.Cells.Item(66,6).Value = 44316
$ws.Cells.Item(66,7).Value = "CALL"
$ws.Cells.Item(66,8).Value = 4175
$ws.Cells.Item(66,9).Value = 20.46
$ws.Cells.Item(66,11).Formula = "=D66*I66"

$ws.Cells.Item(67,1).Value = 44309
$ws.Cells.Item(67,2).Value = "1452-70146"
$ws.Cells.Item(67,3).Value = "BUY"
$ws.Cells.Item(67,4).Value = -100
$ws.Cells.Item(67,5).Value = "SPX"
$ws.Cells.Item(67,6).Value = 44316
$ws.Cells.Item(67,7).Value = "PUT"
$ws.Cells.Item(67,8).Value = 4135
$ws.Cells.Item(67,9).Value = 23.36
$ws.Cells.Item(67,11).Formula = "=D67*I67"
$ws.Cells.Item(67,13).Formula = "=SUM(K64:K67)"

# --- New trade rows 68-71: ID 1452-79806 (26-Apr-2021 expiry SPX straddle) ---
$ws.Cells.Item(68,1).Value = 44309
$ws.Cells.Item(68,2).Value = "1452-79806"
$ws.Cells.Item(68,3).Value = "SELL"
$ws.Cells.Item(68,4).Value = 100
$ws.Cells.Item(68,5).Value = "SPX"
$ws.Cells.Item(68,6).Value = 44312
$ws.Cells.Item(68,7).Value = "CALL"
$ws.Cells.Item(68,8).Value = 4170
$ws.Cells.Item(68,9).Value = 6.83
$ws.Cells.Item(68,11).Formula = "=D68*I68"

$ws.Cells.Item(69,1).Value = 44309
$ws.Cells.Item(69,2).Value = "1452-79806"
$ws.Cells.Item(69,3).Value = "SELL"
$ws.Cells.Item(69,4).Value = 100
$ws.Cells.Item(69,5).Value = "SPX"
$ws.Cells.Item(69,6).Value = 44312
$ws.Cells.Item(69,7).Value = "PUT"
$ws.Cells.Item(69,8).Value = 4140
$ws.Cells.Item(69,9).Value = 8.83
$ws.Cells.Item(69,11).Formula = "=D69*I69"

$ws.Cells.Item(70,1).Value = 44309
$ws.Cells.Item(70,2).Value = "1452-79806"
$ws.Cells.Item(70,3).Value = "BUY"
$ws.Cells.Item(70,4).Value = -100
$ws.Cells.Item(70,5).Value = "SPX"
$ws.Cells.Item(70,6).Value = 44312
$ws.Cells.Item(70,7).Value = "CALL"
$ws.Cells.Item(70,8).Value = 4175
$ws.Cells.Item(70,9).Value = 5.25
$ws.Cells.Item(70,11).Formula = "=D70*I70"

$ws.Cells.Item(71,1).Value = 44309
$ws.Cells.Item(71,2).Value = "1452-79806"
$ws.Cells.Item(71,3).Value = "BUY"
$ws.Cells.Item(71,4).Value = -100
$ws.Cells.Item(71,5).Value = "SPX"
$ws.Cells.Item(71,6).Value = 44312
$ws.Cells.Item(71,7).Value = "PUT"
$ws.Cells.Item(71,8).Value = 3135
$ws.Cells.Item(71,9).Value = 7.51
$ws.Cells.Item(71,11).Formula = "=D71*I71"
$ws.Cells.Item(71,13).Formula = "=SUM(K68:K71)"

# --- Row 72: SDIV stock assignment (ID 1452-89842) ---
$ws.Cells.Item(72,1).Value = 44278
$ws.Cells.Item(72,2).Value = "1452-89842"
$ws.Cells.Item(72,3).Value = "BUY"
$ws.Cells.Item(72,4).Value = -100
$ws.Cells.Item(72,5).Value = "SDIV"
$ws.Cells.Item(72,7).Value = "STOCK"
$ws.Cells.Item(72,9).Value = 14.1
$ws.Cells.Item(72,11).Formula = "=D72*I72"
$ws.Rows.Item(72).Hidden = $true

# --- Row 73: INTC assignment ---
$ws.Cells.Item(73,1).Value = 44309
$ws.Cells.Item(73,2).Value = "1447-49473"
$ws.Cells.Item(73,4).Value = -100
$ws.Cells.Item(73,5).Value = "INTC"
$ws.Cells.Item(73,7).Value = "STOCK"
$ws.Cells.Item(73,8).Value = "ASSIGNMENT"
$ws.Cells.Item(73,9).Value = 61.5
$ws.Cells.Item(73,11).Formula = "=D73*I73"

# --- Mark older short entries (rows 43 & 45) with a '*' note ---
$ws.Cells.Item(43,12).Value = "*"
$ws.Cells.Item(45,12).Value = "*"

# --- Row 74: SPX index assignment ---
$ws.Cells.Item(74,1).Value = 44309
$ws.Cells.Item(74,2).Value = "1444-95611"
$ws.Cells.Item(74,4).Value = -100
$ws.Cells.Item(74,5).Value = "SPX"
$ws.Cells.Item(74,7).Value = "INDEX"
$ws.Cells.Item(74,8).Value = "ASSIGNMENT"
$ws.Cells.Item(74,9).Value = 4170
$ws.Cells.Item(74,11).Formula = "=D74*I74"

# --- Row 75: SPX index exercised ---
$ws.Cells.Item(75,1).Value = 44309
$ws.Cells.Item(75,2).Value = "1444-95611"
$ws.Cells.Item(75,4).Value = 100
$ws.Cells.Item(75,5).Value = "SPX"
$ws.Cells.Item(75,7).Value = "INDEX"
$ws.Cells.Item(75,8).Value = "EXERCISED"
$ws.Cells.Item(75,9).Value = 4175
$ws.Cells.Item(75,11).Formula = "=D75*I75"

# --- Remove the old trailing SUM formula that used to live at M66 ---
$ws.Cells.Item(66,13).ClearContents()

# --- View: reposition active selection to the new last-used cell ---
$ws.Activate()
$ws.Range("L75").Select()

# --- AutoFilter: drop the "this week" dynamic filter, filter Expiry == 23-Apr-2021 ---
$ws.Range("A1:M63").AutoFilter()
$ws.Range("A1:M72").AutoFilter(6, @("23-Apr-2021"), 7)

# --- Re-apply the correct visible/hidden rows (AutoFilter recompute can disturb them) ---
for ($r = 2; $r -le 33; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}
for ($r = 34; $r -le 41; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}
$ws.Rows.Item(42).Hidden = $false
$ws.Rows.Item(43).Hidden = $false
$ws.Rows.Item(44).Hidden = $false
$ws.Rows.Item(45).Hidden = $false
$ws.Rows.Item(46).Hidden = $true
$ws.Rows.Item(47).Hidden = $true
$ws.Rows.Item(48).Hidden = $false
for ($r = 49; $r -le 63; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}
for ($r = 64; $r -le 72; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}
$ws.Rows.Item(73).Hidden = $false
$ws.Rows.Item(74).Hidden = $false
$ws.Rows.Item(75).Hidden = $false
